$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that held "5840751 - Jayne Carlos de Souza Barboza" right under
# "Docentes responsáveis:" (old row 13) is removed; everything below shifts up
# by one row.
$ws.Rows.Item(13).Delete()

# After the shift, a handful of cells need their text content corrected so the
# sheet matches the new data (some values were moved/replaced as part of the
# edit, not simply shifted).

# Row 10 (Objetivos: / was the long "Introduzir..." paragraph) now holds the
# professor identification text.
$ws.Range("B10").Value = "5840751 - Jayne Carlos de Souza Barboza"
$ws.Range("C10").Value = "5840751 - Jayne Carlos de Souza Barboza"

# Row 13 (Programa resumido:) now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) now holds the activation date value.
$ws.Range("B15").Value = "01/01/2020"
$ws.Range("C15").Value = "01/01/2020"

# Row 18 (Método:) now holds the professor identification text.
$ws.Range("B18").Value = "5840751 - Jayne Carlos de Souza Barboza"
$ws.Range("C18").Value = "5840751 - Jayne Carlos de Souza Barboza"

# Row 19 (Critério:) now holds the "Duas provas..." text.
$ws.Range("B19").Value = "Duas provas semestrais teóricas (P1 e P2)."
$ws.Range("C19").Value = "Duas provas semestrais teóricas (P1 e P2)."

# Row 20 (Norma de recuperação:) now holds the average-calculation formula text.
$ws.Range("B20").Value = "A média final (M) será calculada pela expressão M = (P1 + 2 x P2)/3"
$ws.Range("C20").Value = "A média final (M) será calculada pela expressão M = (P1 + 2 x P2)/3"

# Row 21 (Bibliografia:) now holds the recovery-norm text.
$ws.Range("B21").Value = "Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada  recuperação  com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno."
$ws.Range("C21").Value = "Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada  recuperação  com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno."
